# Fruta / hortaliza, semanal
# Insert a new weekly record at row 4 (pushing the existing rows 4-68 down to 5-69)
# and populate the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("4:4").Insert()

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44496
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100106
$ws.Range("H4").Value = "Oleaginosos"
$ws.Range("I4").Value = 100106002
$ws.Range("J4").Value = "Palta"
$ws.Range("K4").Value = "Fuerte"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 62000
$ws.Range("O4").Value = 63000
$ws.Range("P4").Value = 62500
$ws.Range("Q4").Value = "`$/caja 25 kilos"
$ws.Range("R4").Value = "Región de Coquimbo"
$ws.Range("S4").Value = 2500
$ws.Range("T4").Value = 25
